# Happy Farm Fruit liquidation workbook — add the "no vendidos" (unsold)
# items' weight ("2.5") to the main control rows that were still missing it,
# mirroring row 10 (which already carries "2.5" in column E), and leave the
# sheet's selection where the author left it when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rows 11-18 are the other lot/box entries in the main control table; column
# E is "Weight" (重量) and every one of these rows is a 2.5 Kg box, same as
# row 10, but the cells were left blank. Fill them in one shot so they pick
# up the existing shared string "2.5" (same as E10) instead of creating a
# duplicate.
$ws.Range("E11:E18").Value = "2.5"

# Reset the view: scroll back to the top-left (A1) and leave the selection
# on G34, matching where the author's cursor ended up on save.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G34").Select()
